# v2.6 - Added decoupled suspension, four-wheel steering, scripts to
# generate GGV diagram.
#
# Concretely (per the OOXML diff): duplicate the "Truck_Amandla_3Axle"
# sheet into a new "Semi_Truck_Scalable" sheet at the end of the workbook,
# relabel its title cell (H3) to match the new sheet name, make it the
# active sheet/tab, and restore the leftover UI selection state on the
# sheet that used to be active.

$wb = $excel.ActiveWorkbook

# Locate the template sheet to clone and the tail of the sheet list.
$template = $wb.Worksheets.Item("Truck_Amandla_3Axle")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Copy "Truck_Amandla_3Axle" to a new sheet placed after the last sheet.
$template.Copy($null, $lastSheet)

# The newly created sheet is now the last one - rename & relabel it.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Semi_Truck_Scalable"
$newSheet.Range("H3").Value = "Semi_Truck_Scalable"

# Restore the old selection left behind on Truck_Amandla_3Axle (its cursor
# moved on to D24 in the saved workbook).
$template.Activate() | Out-Null
$template.Range("D24").Select() | Out-Null

# The new sheet becomes the active / selected tab, with the cursor resting
# on J17 (bottom-right frozen pane).
$newSheet.Activate() | Out-Null
$newSheet.Range("J17").Select() | Out-Null
